$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 131 (current rows 131:206 shift down to 132:207)
$ws.Rows.Item(131).Insert()

# Populate the new row 131 with the new weekly price record
$ws.Range("A131").Value = 5
$ws.Range("B131").Value = "Macroferia Regional de Talca"
$ws.Range("C131").Value = "Maule"
$ws.Range("D131").Value = 45001
$ws.Range("E131").Value = 7
$ws.Range("F131").Value = 100112030
$ws.Range("G131").Value = "Poroto granado"
$ws.Range("H131").Value = "Sin especificar"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 300
$ws.Range("K131").Value = 28000
$ws.Range("L131").Value = 28000
$ws.Range("M131").Value = 28000
$ws.Range("N131").Value = "$/saco 25 kilos"
$ws.Range("O131").Value = "Región del Maule"
$ws.Range("P131").Value = 1120
$ws.Range("Q131").Value = 25
$ws.Range("R131").Value = "Hortaliza"
